# Udmurt_data.xlsx update
# - Uniformises the valency-class rows on the "Udmurt" sheet: rows that only
#   had column I (X) filled in with "TR" now also get Y (col J) = "ACC" and
#   locus (col K) = "TR", while X itself becomes "NOM". Rows whose X was
#   already "*" just get the same "*" filled into the Y column.
# - Resets the saved view (scroll position / selection) of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Udmurt")

# Rows where X (I) was "TR" with Y (J) and locus (K) empty:
#   I -> NOM, J -> ACC, K -> TR  (L keeps its existing "TR" value)
$trRows = @(5,9,10,14,16,17,19,20,21,27,28,29,32,34,37,40,41,42,44,45,47,50,51,53,56,61,64,67,70,71,72,73,76,77,86,87,92,94,97,98,103,104,106,107,108,109,110,111,116,120,125,127)

foreach ($r in $trRows) {
    $ws.Cells.Item($r, 9).Value = "NOM"
    $ws.Cells.Item($r, 10).Value = "ACC"
    $ws.Cells.Item($r, 11).Value = "TR"
}

# Rows where X (I) was already "*" with Y (J) empty and locus (K) = "*":
#   J -> "*"
$starRows = @(38,46,122,130)

foreach ($r in $starRows) {
    $ws.Cells.Item($r, 10).Value = "*"
}

# Reset the stored view: no scrolled top-left cell, selection back at A1.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
